$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.970.47'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '2.500.26'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').Value = '537.99'
$ws.Range('E5').Value = '  +1.41%  '
$ws.Range('D6').Value = '138.67'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.565'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').Value = '2.525.44'
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('E11').Value = '  -0.04%  '
$ws.Range('E12').Value = '  -1.17%  '
$ws.Range('E13').Value = '  -2.14%  '
$ws.Range('D14').Value = '2.954.73'
$ws.Range('E14').Value = '  -0.04%  '
$ws.Range('D15').Value = '23.25'
$ws.Range('E15').Value = '  +1.03%  '
$ws.Range('D16').Value = '58.870.06'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('D18').Value = '2.526.98'
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('D21').Value = '326.16'
$ws.Range('E21').Value = '  +1.47%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = '5.91'
$ws.Range('E23').Value = '  +1.73%  '
$ws.Range('D24').Value = '65.54'
$ws.Range('E24').Value = '  +5.25%  '
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  +0.46%  '
$ws.Range('D28').Value = '7.67'
$ws.Range('E28').Value = '  -1.54%  '
$ws.Range('D29').Value = '6.78'
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('D30').Value = '0.0₃0780'
$ws.Range('E30').Value = '  +1.46%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').Value = '169.45'
$ws.Range('E32').Value = '  +4.69%  '
$ws.Range('E33').Value = '  +7.50%  '
$ws.Range('E34').Value = '  +3.19%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('E36').Value = '  +0.94%  '
$ws.Range('D37').Value = '4.13'
$ws.Range('E37').Value = '  -1.66%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('D39').Value = '36.66'
$ws.Range('E39').Value = '  -0.84%  '
$ws.Range('D40').Value = '0.830'
$ws.Range('E40').Value = '  +3.62%  '
$ws.Range('D41').Value = '3.65'
$ws.Range('E41').Value = '  +0.50%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value = '5.28'
$ws.Range('E42').Value = '  +1.91%  '
$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').Value = '284.42'
$ws.Range('E43').Value = '  +2.04%  '
$ws.Range('D44').Value = '0.994'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('D45').Value = '130.84'
$ws.Range('E45').Value = '  +7.65%  '
$ws.Range('D46').Value = '0.605'
$ws.Range('E46').Value = '  +1.92%  '
$ws.Range('D47').Value = '10.89'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('D49').Value = '0.0513'
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('D50').Value = '0.0223'
$ws.Range('E50').Value = '  +0.39%  '
$ws.Range('D51').Value = '17.57'
$ws.Range('E51').Value = '  +0.33%  '
